$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.019999999999999
$ws.Range("C2").Value = 1.054472091178195
$ws.Range("D2").Value = 1.052062178747892
$ws.Range("E2").Value = 1.060214665458465
$ws.Range("F2").Value = 1.069049703663304
$ws.Range("I2").Value = 1.044997968271764
$ws.Range("J2").Value = 1.059484010098365
$ws.Range("K2").Value = 1.05481195174809
$ws.Range("L2").Value = 1.062942053370013
$ws.Range("M2").Value = 1.071753249373306
$ws.Range("N2").Value = 1.005712725503983
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.055658446177062
$ws.Range("D3").Value = 1.052942796986936
$ws.Range("E3").Value = 1.061308740310813
$ws.Range("F3").Value = 1.070313874540435
$ws.Range("I3").Value = 1.045327820263063
$ws.Range("J3").Value = 1.060320512765051
$ws.Range("K3").Value = 1.055505166449624
$ws.Range("L3").Value = 1.063849803412116
$ws.Range("M3").Value = 1.072832404193492
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.056425902941882
$ws.Range("D4").Value = 1.05351237266286
$ws.Range("E4").Value = 1.062016841348913
$ws.Range("F4").Value = 1.07113235149337
$ws.Range("I4").Value = 1.045539949128471
$ws.Range("J4").Value = 1.060861044696931
$ws.Range("K4").Value = 1.055952842410645
$ws.Range("L4").Value = 1.06443674494024
$ws.Range("M4").Value = 1.073530600518336
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.056748496868093
$ws.Range("D5").Value = 1.053751764486733
$ws.Range("E5").Value = 1.062314566614622
$ws.Range("F5").Value = 1.071476553976797
$ws.Range("I5").Value = 1.045628815433396
$ws.Range("J5").Value = 1.061088107726274
$ws.Range("K5").Value = 1.05614083521901
$ws.Range("L5").Value = 1.064683392279685
$ws.Range("M5").Value = 1.07382410207484
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.056802659221685
$ws.Range("D6").Value = 1.053791956051097
$ws.Range("E6").Value = 1.0623645583456
$ws.Range("F6").Value = 1.071534353857305
$ws.Range("I6").Value = 1.045643718151541
$ws.Range("J6").Value = 1.061126222277336
$ws.Range("K6").Value = 1.056172387731416
$ws.Range("L6").Value = 1.064724799446275
$ws.Range("M6").Value = 1.073873381127233
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.056430213632117
$ws.Range("D7").Value = 1.053515571654665
$ws.Range("E7").Value = 1.062020819411753
$ws.Range("F7").Value = 1.071136950291158
$ws.Range("I7").Value = 1.045541137792552
$ws.Range("J7").Value = 1.060864079417071
$ws.Range("K7").Value = 1.055955355204836
$ws.Range("L7").Value = 1.064440041057196
$ws.Range("M7").Value = 1.073534522378374
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.054873066591359
$ws.Range("D8").Value = 1.052359838425651
$ws.Range("E8").Value = 1.060584380312811
$ws.Range("F8").Value = 1.069476838813295
$ws.Range("I8").Value = 1.045109714207893
$ws.Range("J8").Value = 1.05976686408343
$ws.Range("K8").Value = 1.055046409525277
$ws.Range("L8").Value = 1.063248922030867
$ws.Range("M8").Value = 1.072117973906296
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.05212761472077
$ws.Range("D9").Value = 1.050321406101122
$ws.Range("E9").Value = 1.058054385209999
$ws.Range("F9").Value = 1.066555079713012
$ws.Range("I9").Value = 1.044339462350879
$ws.Range("J9").Value = 1.057827714837714
$ws.Range("K9").Value = 1.053437962622458
$ws.Range("L9").Value = 1.061146658248612
$ws.Range("M9").Value = 1.069621101950413
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.050296161562438
$ws.Range("D10").Value = 1.048961154695234
$ws.Range("E10").Value = 1.056368462766407
$ws.Range("F10").Value = 1.06460955401516
$ws.Range("I10").Value = 1.043819196306215
$ws.Range("J10").Value = 1.056531049432969
$ws.Range("K10").Value = 1.052361074897677
$ws.Range("L10").Value = 1.059742831601926
$ws.Range("M10").Value = 1.067955954523955
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.049502825210934
$ws.Range("D11").Value = 1.048371835619489
$ws.Range("E11").Value = 1.055638600882923
$ws.Range("F11").Value = 1.06376764870004
$ws.Range("I11").Value = 1.043592305466183
$ws.Range("J11").Value = 1.05596864089425
$ws.Range("K11").Value = 1.051893672619785
$ws.Range("L11").Value = 1.05913439575592
$ws.Range("M11").Value = 1.067234775882555
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.049208096796691
$ws.Range("D12").Value = 1.048152887250181
$ws.Range("E12").Value = 1.055367519411918
$ws.Range("F12").Value = 1.063455003907305
$ws.Range("I12").Value = 1.043507785362173
$ws.Range("J12").Value = 1.055759594280528
$ws.Range("K12").Value = 1.051719892024483
$ws.Range("L12").Value = 1.058908308732193
$ws.Range("M12").Value = 1.066966872589969
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.04927131924716
$ws.Range("D13").Value = 1.048199854649414
$ws.Range("E13").Value = 1.055425666346157
$ws.Range("F13").Value = 1.063522063851931
$ws.Range("I13").Value = 1.043525926204532
$ws.Range("J13").Value = 1.055804441981598
$ws.Range("K13").Value = 1.051757176104533
$ws.Range("L13").Value = 1.058956809117386
$ws.Range("M13").Value = 1.067024339931927
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.049478463859229
$ws.Range("D14").Value = 1.048353738283502
$ws.Range("E14").Value = 1.055616192768541
$ws.Range("F14").Value = 1.063741803837424
$ws.Range("I14").Value = 1.043585323962968
$ws.Range("J14").Value = 1.055951363961104
$ws.Range("K14").Value = 1.051879311264907
$ws.Range("L14").Value = 1.059115709111322
$ws.Range("M14").Value = 1.067212631440789
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.049606086070821
$ws.Range("D15").Value = 1.048448544550085
$ws.Range("E15").Value = 1.05573358522431
$ws.Range("F15").Value = 1.063877202954436
$ws.Range("I15").Value = 1.043621888700637
$ws.Range("J15").Value = 1.056041868450339
$ws.Range("K15").Value = 1.051954540668083
$ws.Range("L15").Value = 1.059213601093919
$ws.Range("M15").Value = 1.067328640636334
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.050348806073681
$ws.Range("D16").Value = 1.049000259028266
$ws.Range("E16").Value = 1.056416904417637
$ws.Range("F16").Value = 1.064665439301899
$ws.Range("I16").Value = 1.043834220291604
$ws.Range("J16").Value = 1.056568354643854
$ws.Range("K16").Value = 1.052392071536554
$ws.Range("L16").Value = 1.059783199374473
$ws.Range("M16").Value = 1.068003813294832
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.050814611529737
$ws.Range("D17").Value = 1.049346248402561
$ws.Range("E17").Value = 1.056845572677898
$ws.Range("F17").Value = 1.065160017113709
$ws.Range("I17").Value = 1.043966978282587
$ws.Range("J17").Value = 1.056898351953081
$ws.Range("K17").Value = 1.052666227211551
$ws.Range("L17").Value = 1.060140339916892
$ws.Range("M17").Value = 1.068427287635777
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.051086278682428
$ws.Range("D18").Value = 1.049548027207231
$ws.Range("E18").Value = 1.05709562253874
$ws.Range("F18").Value = 1.065448546292963
$ws.Range("I18").Value = 1.044044258253205
$ws.Range("J18").Value = 1.057090742814562
$ws.Range("K18").Value = 1.052826031235819
$ws.Range("L18").Value = 1.060348599076377
$ws.Range("M18").Value = 1.068674277836398
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.051178905258109
$ws.Range("D19").Value = 1.049616823339727
$ws.Range("E19").Value = 1.057180885650816
$ws.Range("F19").Value = 1.065546935875442
$ws.Range("I19").Value = 1.044070582363337
$ws.Range("J19").Value = 1.057156327731295
$ws.Range("K19").Value = 1.052880502268081
$ws.Range("L19").Value = 1.060419600798656
$ws.Range("M19").Value = 1.068758492578642
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.050764638057319
$ws.Range("D20").Value = 1.049309130223092
$ws.Range("E20").Value = 1.056799579097072
$ws.Range("F20").Value = 1.065106948404302
$ws.Range("I20").Value = 1.043952750696383
$ws.Range("J20").Value = 1.056862955784694
$ws.Range("K20").Value = 1.052636823908037
$ws.Range("L20").Value = 1.060102027801797
$ws.Range("M20").Value = 1.068381854431901
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.049417466275324
$ws.Range("D21").Value = 1.048308424748189
$ws.Range("E21").Value = 1.055560086896958
$ws.Range("F21").Value = 1.063677093796196
$ws.Range("I21").Value = 1.043567839508711
$ws.Range("J21").Value = 1.055908103058655
$ws.Range("K21").Value = 1.051843350109905
$ws.Range("L21").Value = 1.059068919431616
$ws.Range("M21").Value = 1.067157184999248
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.048570167306032
$ws.Range("D22").Value = 1.047678957473418
$ws.Range("E22").Value = 1.054780894002622
$ws.Range("F22").Value = 1.062778528155716
$ws.Range("I22").Value = 1.043324425866371
$ws.Range("J22").Value = 1.055306921045502
$ws.Range("K22").Value = 1.051343497424134
$ws.Range("L22").Value = 1.058418860016444
$ws.Range("M22").Value = 1.066387038781507
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.04901936364468
$ws.Range("D23").Value = 1.048012677163748
$ws.Range("E23").Value = 1.055193947499533
$ws.Range("F23").Value = 1.063254833638198
$ws.Range("I23").Value = 1.043453597369475
$ws.Range("J23").Value = 1.055625697923545
$ws.Range("K23").Value = 1.051608570487221
$ws.Range("L23").Value = 1.058763516934851
$ws.Range("M23").Value = 1.066795322412528
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.050787219008352
$ws.Range("D24").Value = 1.04932590242599
$ws.Range("E24").Value = 1.056820361569978
$ws.Range("F24").Value = 1.065130927713813
$ws.Range("I24").Value = 1.043959180010522
$ws.Range("J24").Value = 1.056878950070336
$ws.Range("K24").Value = 1.052650110323532
$ws.Range("L24").Value = 1.060119339566659
$ws.Range("M24").Value = 1.068402383787652
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.052837575302802
$ws.Range("D25").Value = 1.050848615973955
$ws.Range("E25").Value = 1.058708314441935
$ws.Range("F25").Value = 1.067310010345689
$ws.Range("I25").Value = 1.044539781627034
$ws.Range("J25").Value = 1.058329714527376
$ws.Range("K25").Value = 1.053854590812245
$ws.Range("L25").Value = 1.061690548001802
$ws.Range("M25").Value = 1.070266697583149
